# Edit script: add "SecondSheet" worksheet, insert a new parameter row into
# ValidSheet, and populate SecondSheet with sample parameter data.

$wb = $excel.ActiveWorkbook

# --- Sheet references (before any insertion) ---
$validSheet = $wb.Worksheets.Item(1)          # "ValidSheet"

# --- 1. Insert a new row into ValidSheet (row 2) ------------------------
$validSheet.Rows.Item(2).Insert()
$validSheet.Range("A2").Value = "Path1"
$validSheet.Range("B2").Value = "Param1"
$validSheet.Range("C2").Value = 0

# --- 2. Add the new "SecondSheet" worksheet right after "ValidSheet" ----
$secondSheet = $wb.Worksheets.Add($null, $validSheet)
$secondSheet.Name = "SecondSheet"

# Header row
$secondSheet.Range("A1").Value = "Container Path"
$secondSheet.Range("B1").Value = "Parameter Name"
$secondSheet.Range("C1").Value = "Value"
$secondSheet.Range("D1").Value = "Units"

# Data rows
$secondSheet.Range("A2").Value = "Path1"
$secondSheet.Range("B2").Value = "Param1"
$secondSheet.Range("C2").Value = 5
$secondSheet.Range("D2").Value = "mg"

$secondSheet.Range("A3").Value = "Path2"
$secondSheet.Range("B3").Value = "DistincParam"
$secondSheet.Range("C3").Value = 1
$secondSheet.Range("D3").Value = "µmol"

$secondSheet.Range("A4").Value = "Applications|Glucose_iv_infusion"
$secondSheet.Range("B4").Value = "Active"
$secondSheet.Range("C4").Value = 0

$secondSheet.Columns.Item(1).ColumnWidth = 30.08984375
$secondSheet.Columns.Item(2).ColumnWidth = 15.36328125
$secondSheet.Columns.Item(3).ColumnWidth = 5.81640625

# --- 3. Selections (mirrors final state of the workbook) ----------------
$validSheet.Select()
$validSheet.Rows.Item(3).Select()

$secondSheet.Select()
$secondSheet.Range("C5").Select()
